$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '27.705.54'
$ws.Range('E2').Value = '  -1.92%  '

$ws.Range('D3').Value = '1.757.05'
$ws.Range('E3').Value = '  -2.17%  '

$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '1.005'
$ws.Range('E4').Value = '  +0.12%  '

$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '324.92'
$ws.Range('E5').Value = '  -4.00%  '

$ws.Range('E6').Value = '  +0.09%  '

$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.4497'
$ws.Range('E7').Value = '  -1.58%  '

$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.3739'
$ws.Range('E8').Value = '  -0.11%  '

$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '45.40'
$ws.Range('E9').Value = '  +0.60%  '

$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.07552'
$ws.Range('E10').Value = '  -0.49%  '

$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '1.126'
$ws.Range('E11').Value = '  -1.57%  '

$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '1.002'
$ws.Range('E12').Value = '  -0.02%  '

$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '21.76'
$ws.Range('E13').Value = '  -2.41%  '

$ws.Range('E14').Value = '  -1.22%  '

$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '7.354'
$ws.Range('E15').Value = '  -1.76%  '

$ws.Range('D16').Value = '1.762.51'
$ws.Range('E16').Value = '  -2.00%  '

$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '0.00001075'
$ws.Range('E17').Value = '  -1.22%  '

$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '87.87'
$ws.Range('E18').Value = '  +8.28%  '

$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '0.06221'
$ws.Range('E19').Value = '  -7.64%  '

$ws.Range('E20').Value = '  +0.05%  '

$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '17.27'
$ws.Range('E21').Value = '  -0.79%  '

$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '6.183'
$ws.Range('E22').Value = '  -2.92%  '

$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '0.5325'
$ws.Range('E23').Value = '  -3.28%  '

$ws.Range('D24').Value = '27.743.42'
$ws.Range('E24').Value = '  -1.81%  '

$ws.Range('E25').Value = '  -1.53%  '

$ws.Range('E26').Value = '  -4.28%  '

$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '20.67'
$ws.Range('E27').Value = '  +0.06%  '

$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '152.96'
$ws.Range('E28').Value = '  +0.66%  '

$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '2.362'
$ws.Range('E29').Value = '  +0.54%  '

$ws.Range('D30').Value = '1.955.69'
$ws.Range('E30').Value = '  -2.36%  '

$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '128.52'
$ws.Range('E31').Value = '  -3.10%  '

$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '1.217'
$ws.Range('E32').Value = '  -1.20%  '

$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '0.09328'
$ws.Range('E33').Value = '  -1.70%  '

$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '5.737'
$ws.Range('E34').Value = '  -1.07%  '

$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '3.640'
$ws.Range('E35').Value = '  -9.68%  '

$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '12.66'
$ws.Range('E36').Value = '  +5.33%  '

$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.02332'
$ws.Range('E37').Value = '  -0.41%  '

$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.2169'
$ws.Range('E38').Value = '  -7.81%  '

$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.06149'
$ws.Range('E39').Value = '  -2.74%  '

$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.6482'
$ws.Range('E40').Value = '  -1.64%  '

$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '5.080'
$ws.Range('E41').Value = '  -3.24%  '

$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '1.200'
$ws.Range('E42').Value = '  -1.83%  '

$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '7.998'
$ws.Range('E43').Value = '  -4.23%  '

$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '1.420'
$ws.Range('E44').Value = '  -4.34%  '

$ws.Range('E45').Value = '  +0.08%  '

$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '13.82'
$ws.Range('E46').Value = '  -3.10%  '

$ws.Range('E47').Value = '  -1.63%  '

$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '3.752'
$ws.Range('E48').Value = '  -2.61%  '

$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '126.27'
$ws.Range('E49').Value = '  -3.18%  '

$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '1.989'
$ws.Range('E50').Value = '  -1.97%  '

$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.06907'
$ws.Range('E51').Value = '  -3.31%  '
